$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "2024-09-29 00:05:44"
$ws.Range("B4").Value = "MOCK_check_availability"
$ws.Range("C4").Value = "MOCKURL_https://www.opentable.com/r/bar-spero-washington/"
$ws.Range("D4").Value = "MOCK_No availability for the selected date."

# "2024-09-29" alone would be auto-recognized as a date by Excel's type
# inference, so force it to stay plain text with a leading apostrophe and
# then clear the resulting formatting so no style index is left behind.
$ws.Range("E4").Value = "'2024-09-29"
$ws.Range("E4").ClearFormats()

$ws.Range("F4").Value = "00:05:44"
